# clean the execute data
# Clear the "实际结果" (actual result, column F) and "是否通过" (pass/fail, column G)
# columns for the data rows (2-9), leaving the "预期结果" (expected result,
# column E) values untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove computed/actual-result and pass/fail values for every data row.
$ws.Range("F2:G9").ClearContents()
